# Natmi following Dr Hou advice
#
# Update the LR-pairs sheet: the two existing communication edges (row 2:
# FAPs -> sCs, row 3: sCs -> sCs [was mis-scaled]) are refreshed with new
# numbers, and a third sending cluster "ECs" is added as a new row, each
# reporting the same Il18 -> Il1rapl1 ligand/receptor pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: ECs -> Il18 -> Il1rapl1 -> sCs ----
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Il18"
$ws.Range("C2").Value2 = "Il1rapl1"
$ws.Range("D2").Value2 = "sCs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 2.878986666666667
$ws.Range("H2").Value2 = 8.63696
$ws.Range("I2").Value2 = 0.2331567682967092
$ws.Range("J2").Value2 = 0.2331567682967092
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.3962093333333334
$ws.Range("N2").Value2 = 1.188628
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 1.140681387875556
$ws.Range("R2").Value2 = 10.26613249088
$ws.Range("S2").Value2 = 0.2331567682967092
$ws.Range("T2").Value2 = 0.2331567682967092

# ---- Row 3: FAPs -> Il18 -> Il1rapl1 -> sCs ----
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Il18"
$ws.Range("C3").Value2 = "Il1rapl1"
$ws.Range("D3").Value2 = "sCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 8.325812333333333
$ws.Range("H3").Value2 = 24.977437
$ws.Range("I3").Value2 = 0.6742717913773656
$ws.Range("J3").Value2 = 0.6742717913773655
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.3962093333333334
$ws.Range("N3").Value2 = 1.188628
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 3.298764554048445
$ws.Range("R3").Value2 = 29.688880986436
$ws.Range("S3").Value2 = 0.6742717913773656
$ws.Range("T3").Value2 = 0.6742717913773655

# ---- Row 4 (new): sCs -> Il18 -> Il1rapl1 -> sCs ----
$ws.Range("A4").Value2 = "sCs"
$ws.Range("B4").Value2 = "Il18"
$ws.Range("C4").Value2 = "Il1rapl1"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 1.143059
$ws.Range("H4").Value2 = 3.429177
$ws.Range("I4").Value2 = 0.09257144032592537
$ws.Range("J4").Value2 = 0.09257144032592536
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.3962093333333334
$ws.Range("N4").Value2 = 1.188628
$ws.Range("O4").Value2 = 1
$ws.Range("P4").Value2 = 1
$ws.Range("Q4").Value2 = 0.4528906443506667
$ws.Range("R4").Value2 = 4.076015799156
$ws.Range("S4").Value2 = 0.09257144032592537
$ws.Range("T4").Value2 = 0.09257144032592536
